# Update gh-pages output data (苏州-漫展信息.xlsx)
# Refresh "想去人数" (F) counters and flip the "不可售" (not-for-sale)
# marker in G to its numeric minimum-price value on both the "展览"
# and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F5").Value = 11323
        $ws.Range("F9").Value = 11264
        $ws.Range("F11").Value = 1149
        $ws.Range("F12").Value = 64
        $ws.Range("F13").Value = 1738
        $ws.Range("F14").Value = 5623
        $ws.Range("G18").Value = 58
    }
    else {
        $ws.Range("F7").Value = 11323
        $ws.Range("F11").Value = 11264
        $ws.Range("F13").Value = 1149
        $ws.Range("F14").Value = 64
        $ws.Range("F15").Value = 1738
        $ws.Range("F16").Value = 5623
        $ws.Range("G20").Value = 58
    }
}
